# Generate Report for Handback
# A new file (21231fb8-99ee-4e9a-b5da-0b547ff5db4a) has been handed back.
# It becomes the new "row 2" entry on every sheet, and the previously
# reported file (82d7a099-8572-426d-b561-4a6d0d558207) is pushed down to
# "row 3" (same data it already had).

$wb = $excel.ActiveWorkbook

$dateFmt = "yyyy-mm-dd HH:mm:ss"

$oldGuid = "82d7a099-8572-426d-b561-4a6d0d558207"
$newGuid = "21231fb8-99ee-4e9a-b5da-0b547ff5db4a"

$oldMd = "$oldGuid.md"
$newMd = "$newGuid.md"

# ---------------------------------------------------------------------
# Overview sheet: File Name | zh-cn | de-de
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

# Clear existing hyperlinks so they get re-created (and re-numbered) below.
$ws.Range("A2").Hyperlinks.Delete()

# Row 2 -> the newly handed-back file.
$ws.Range("A2").Value = $newMd
$ws.Range("B2").Value = "Handed back: in sync with en-US"
$ws.Range("C2").Value = "Handed back: in sync with en-US"

# Row 3 -> the file that used to be reported in row 2.
$ws.Range("A3").Value = $oldMd
$ws.Range("B3").Value = "Handed back: in sync with en-US"
$ws.Range("C3").Value = "Handed back: in sync with en-US"

$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/e8893fb22cbc1b958512083c8bd234abb2f1b86c/e2e/$newMd", "", "", $newMd) | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/e8893fb22cbc1b958512083c8bd234abb2f1b86c/e2e/$oldMd", "", "", $oldMd) | Out-Null

# ---------------------------------------------------------------------
# Per-locale detail sheets (zh-cn / de-de)
# ---------------------------------------------------------------------
$locales = @(
    @{ Sheet = "zh-cn"; Code = "zh-cn";
       NewXlf = "$newGuid.cbff2a9cde1128f13779c41f25d3a7390c88e056.zh-cn.xlf";
       OldXlf = "$oldGuid.1a5694b0db3e42df80cf4ac3e29058cff8eda798.zh-cn.xlf";
       NewHandoff = "2016-03-20 04:34:32"; NewHandback = "2016-03-20 04:35:13";
       OldHandoff = "2016-03-20 04:31:10"; OldHandback = "2016-03-20 04:32:14";
       HandoffOrg = "oltest-zhcn-fly"; HandoffSha = "7cc607d7f95ef3ad8648a099792d55709842ac18"; HandoffFileSha = "46160d8d80b20e2feaffe35110bbf74af524f71f";
       HandbackSha = "009cca8426b314f264284b735c3aa4e83b632f2b" },
    @{ Sheet = "de-de"; Code = "de-de";
       NewXlf = "$newGuid.cbff2a9cde1128f13779c41f25d3a7390c88e056.de-de.xlf";
       OldXlf = "$oldGuid.1a5694b0db3e42df80cf4ac3e29058cff8eda798.de-de.xlf";
       NewHandoff = "2016-03-20 04:34:41"; NewHandback = "2016-03-20 04:35:28";
       OldHandoff = "2016-03-20 04:31:19"; OldHandback = "2016-03-20 04:32:29";
       HandoffOrg = "oltest-dede-fly"; HandoffSha = "8de731932493eaf71c47ebbe7b3e1439f8f68bd8"; HandoffFileSha = "a29d03243dde5967aa6cbb7c47328d97aa339ad3";
       HandbackSha = "b495b1e95dc91e4318e9ef6b765cc59857939cbf" }
)

foreach ($loc in $locales) {
    $ws = $wb.Worksheets.Item($loc.Sheet)

    # Remove the old hyperlinks tied to row 2 so fresh ones (pointing at the
    # new file first, then the old file) get assigned in the right order.
    $ws.Range("A2").Hyperlinks.Delete()
    $ws.Range("D2").Hyperlinks.Delete()
    $ws.Range("F2").Hyperlinks.Delete()
    $ws.Range("G2").Hyperlinks.Delete()

    # --- Row 2: newly handed-back file -------------------------------
    $ws.Range("A2").Value = $newMd
    $ws.Range("B2").Value = ".md"
    $ws.Range("C2").Value = "Handed back: in sync with en-US"
    $ws.Range("D2").Value = $loc.NewXlf
    $ws.Range("E2").Value = $loc.NewHandoff
    $ws.Range("E2").NumberFormat = $dateFmt
    $ws.Range("F2").Value = $newMd
    $ws.Range("G2").Value = $loc.NewXlf
    $ws.Range("H2").Value = $loc.NewHandback
    $ws.Range("H2").NumberFormat = $dateFmt
    $ws.Range("J2").Value = "Include"

    # --- Row 3: file previously reported in row 2 ---------------------
    $ws.Range("A3").Value = $oldMd
    $ws.Range("B3").Value = ".md"
    $ws.Range("C3").Value = "Handed back: in sync with en-US"
    $ws.Range("D3").Value = $loc.OldXlf
    $ws.Range("E3").Value = $loc.OldHandoff
    $ws.Range("E3").NumberFormat = $dateFmt
    $ws.Range("F3").Value = $oldMd
    $ws.Range("G3").Value = $loc.OldXlf
    $ws.Range("H3").Value = $loc.OldHandback
    $ws.Range("H3").NumberFormat = $dateFmt
    $ws.Range("J3").Value = "Include"

    $mdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/e8893fb22cbc1b958512083c8bd234abb2f1b86c/e2e"
    $handoffBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$($loc.HandoffSha)/ol-handoff/OpenLocalizationTestOrg/$($loc.HandoffOrg)/yuwzho/ht"
    $handoffMdBase = "https://github.com/OpenLocalizationTestOrg/$($loc.HandoffOrg)/blob/$($loc.HandoffFileSha)/e2e"
    $handbackBase = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/$($loc.HandbackSha)/ol-handback/OpenLocalizationTestOrg/$($loc.HandoffOrg)/yuwzho/ht"

    # New-file hyperlinks first -> get rId2..rId5
    $ws.Hyperlinks.Add($ws.Range("A2"), "$mdUrl/$newMd", "", "", $newMd) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("D2"), "$handoffBase/$($loc.NewXlf)", "", "", $loc.NewXlf) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("F2"), "$handoffMdBase/$newMd", "", "", $newMd) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("G2"), "$handbackBase/$($loc.NewXlf)", "", "", $loc.NewXlf) | Out-Null

    # Old-file hyperlinks second -> get rId6..rId9
    $ws.Hyperlinks.Add($ws.Range("A3"), "$mdUrl/$oldMd", "", "", $oldMd) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("D3"), "$handoffBase/$($loc.OldXlf)", "", "", $loc.OldXlf) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("F3"), "$handoffMdBase/$oldMd", "", "", $oldMd) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("G3"), "$handbackBase/$($loc.OldXlf)", "", "", $loc.OldXlf) | Out-Null
}

Write-Output "Handback report updated"
